$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "whed_inst": insert an "Address Concat" column (new col G),
# and add a sample data row (row 2) with a TEXTJOIN formula.
# ------------------------------------------------------------------
$wsInst = $wb.Worksheets.Item("whed_inst")

$wsInst.Range("G1").EntireColumn.Insert()
$wsInst.Range("G1").Value = "Address Concat"

$wsInst.Range("H2").Value = "Unit 1, 1 Smith Street, Pall Mall, South Kensington"
$wsInst.Range("I2").Value = "Kensington"
$wsInst.Range("J2").Value = "SW11"
$wsInst.Range("K2").Value = "Michigan"
$wsInst.Range("M2").Value = "Australia"
$wsInst.Range("G2").Formula = "=TEXTJOIN("" "", TRUE, H2:M2)"

$wsInst.Range("G2").Select() | Out-Null

# ------------------------------------------------------------------
# Sheet "ext_inst": insert a "Row Index" column (new col A), rename
# the former "Institution Type" column to "Postal Address Country",
# and add a sample data row (row 2) with a TEXTJOIN formula.
# ------------------------------------------------------------------
$wsExtInst = $wb.Worksheets.Item("ext_inst")

$wsExtInst.Range("A1").EntireColumn.Insert()
$wsExtInst.Range("A1").Value = "Row Index"
$wsExtInst.Range("A1").Font.Bold = $true

$wsExtInst.Range("N1").Value = "Postal Address Country"

$wsExtInst.Range("A2").Value = 1
$wsExtInst.Range("B2").Value = "Blue62"
$wsExtInst.Range("C2").Value = "Acorn Institute"
$wsExtInst.Range("D2").Value = "Oak Academy"
$wsExtInst.Range("G2").Value = "Unit 1"
$wsExtInst.Range("H2").Value = "1 Smith Street"
$wsExtInst.Range("I2").Value = "Pall Mall"
$wsExtInst.Range("J2").Value = "South Kensington"
$wsExtInst.Range("K2").Value = "Kensington"
$wsExtInst.Range("L2").Value = "Michigan"
$wsExtInst.Range("M2").Value = "SW11"
$wsExtInst.Range("F2").Formula = "=TEXTJOIN("" "", TRUE, G2:M2)"

$wsExtInst.Range("O1").Select() | Out-Null
$wsExtInst.Activate() | Out-Null

$wb.Save() | Out-Null
